$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B ("Name") to hold the new "Skill_ID" field.
#    This shifts Name/DamageCoefficient/CoolDown/UseRange/SkillPriority one column to the right
#    (B->C, C->D, D->E, E->F, F->G), carrying along styles/values/widths automatically.
$ws.Columns("B:B").Insert()

# 2. Populate the new Skill_ID column (header / type / label / data rows).
$ws.Range("B1").Value = "Skill_ID"
$ws.Range("B2").Value = "string"
$ws.Range("B3").Value = "스킬식별자"
$ws.Range("B4").Value = "SKL_SWIP"
$ws.Range("B5").Value = "SKL_HEADBUTT"
$ws.Range("B6").Value = "SKL_DROPKICK"

# 3. Remove the old "SkillPriority" column, which is now column G after the insert above.
$ws.Columns("G:G").Delete()

# 4. Rename the skill entries in column C (previously column B).
$ws.Range("C4").Value = "Swip"
$ws.Range("C5").Value = "HeadButt"
$ws.Range("C6").Value = "DropKick"

# 5. Update the numeric stat columns (D: DamageCoefficient, E: CoolDown, F: UseRange).
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.5
$ws.Range("F4").Value = 1

$ws.Range("D5").Value = 1.3
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.5

$ws.Range("D6").Value = 1.6
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 3.5

# 6. Adjust the new Skill_ID column width; the remaining columns already carried
#    their previous widths across from the insert operation.
$ws.Columns("B:B").ColumnWidth = 14.5

# 7. Update the active selection to match the saved state.
$ws.Range("F6").Select()

# 8. Restore the page setup (A4, portrait) used for printing this sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
